$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update monetary values to include thousands separator (comma) formatting.
# These are stored as literal text strings (e.g. "2000.00 USD" -> "2,000.00 USD").
$ws.Range("B28").Value = "2,000.00 USD"

$ws.Range("C29").Value = "1,500.00 USD"
$ws.Range("D29").Value = "1,500.00 USD"

$ws.Range("B30").Value = "2,000.00 USD"
$ws.Range("C30").Value = "2,000.00 USD"
$ws.Range("D30").Value = "2,000.00 USD"

$ws.Range("B31").Value = "2,000.00 USD"
$ws.Range("C31").Value = "2,000.00 USD"
$ws.Range("D31").Value = "2,000.00 USD"

$ws.Range("B33").Value = "2,000.00 USD"

$ws.Range("C39").Value = "1,500.00 USD"
$ws.Range("D39").Value = "1,500.00 USD"

# Remove the two trailing blank rows (49 and 50) that held no data.
$ws.Rows("49:50").Delete()

# Restore the view/selection state recorded for this sheet.
$ws.Range("A19").Select()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("E33").Select()
